$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 740
$ws.Range("I41").Value = 240
$ws.Range("J41").Value = 954.2857
$ws.Range("K41").Value = 240
$ws.Range("L41").Value = 954.2857
$ws.Range("M41").Value = 200
$ws.Range("N41").Value = -1834.2857

$ws.Range("H113").Value = 2655.6897
$ws.Range("I113").Value = 3001.0715
$ws.Range("K113").Value = 3001.0715
$ws.Range("M113").Value = 252.9285

$ws.Range("H116").Value = 4081.1904
$ws.Range("I116").Value = 4021.3572
$ws.Range("K116").Value = 4021.3572
$ws.Range("M116").Value = -579.3571999999999

$ws.Range("H132").Value = 2314.6292
$ws.Range("I132").Value = 2232.776
$ws.Range("K132").Value = 6698.328
$ws.Range("M132").Value = -4168.328

$ws.Range("H135").Value = 876.6326
$ws.Range("I135").Value = 525.725
$ws.Range("J135").Value = 2436.2222
$ws.Range("K135").Value = 4731.525000000001
$ws.Range("L135").Value = 21925.9998
$ws.Range("M135").Value = -2196.525000000001
$ws.Range("N135").Value = -26995.9998

$ws.Range("H137").Value = 1140
$ws.Range("I137").Value = 1140
$ws.Range("K137").Value = 3420
$ws.Range("M137").Value = -870

$ws.Range("H138").Value = 1191.0303
$ws.Range("I138").Value = 616.07574
$ws.Range("J138").Value = 2340.9395
$ws.Range("K138").Value = 1848.22722
$ws.Range("L138").Value = 7022.818499999999
$ws.Range("M138").Value = 3291.77278
$ws.Range("N138").Value = -17302.8185

$ws.Range("H141").Value = 2464.2456
$ws.Range("I141").Value = 596.0769
$ws.Range("J141").Value = 6511.9443
$ws.Range("K141").Value = 1788.2307
$ws.Range("L141").Value = 19535.8329
$ws.Range("M141").Value = 3391.7693
$ws.Range("N141").Value = -29895.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3639352.2
$ws.Range("I32").Value = 4584681
$ws.Range("J32").Value = 7300.1577
$ws.Range("K32").Value = 4584681
$ws.Range("L32").Value = 7300.1577
$ws.Range("M32").Value = -4584394
$ws.Range("N32").Value = -7874.1577

$ws.Range("H61").Value = 1126.9412
$ws.Range("I61").Value = 844.24
$ws.Range("J61").Value = 1912.2222
$ws.Range("K61").Value = 844.24
$ws.Range("L61").Value = 1912.2222
$ws.Range("M61").Value = -632.24
$ws.Range("N61").Value = -2336.2222

$ws.Range("H74").Value = 1034.4722
$ws.Range("I74").Value = 1010.0968
$ws.Range("K74").Value = 1010.0968
$ws.Range("M74").Value = -136.0968

$ws.Range("H77").Value = 1034.4722
$ws.Range("I77").Value = 1010.0968
$ws.Range("K77").Value = 5050.484
$ws.Range("M77").Value = -682.4840000000004

$ws.Range("H110").Value = 603.7646999999999
$ws.Range("I110").Value = 561.8333
$ws.Range("J110").Value = 704.4
$ws.Range("K110").Value = 561.8333
$ws.Range("L110").Value = 704.4
$ws.Range("M110").Value = 1483.1667
$ws.Range("N110").Value = -4794.4

$ws.Range("H132").Value = 1777.1632
$ws.Range("I132").Value = 1095.4828
$ws.Range("J132").Value = 2765.6
$ws.Range("K132").Value = 3286.4484
$ws.Range("L132").Value = 8296.799999999999
$ws.Range("M132").Value = -756.4484000000002
$ws.Range("N132").Value = -13356.8

$ws.Range("H136").Value = 1126.9412
$ws.Range("I136").Value = 844.24
$ws.Range("J136").Value = 1912.2222
$ws.Range("K136").Value = 2532.72
$ws.Range("L136").Value = 5736.6666
$ws.Range("M136").Value = 17.27999999999975
$ws.Range("N136").Value = -10836.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 45456264
$ws.Range("I99").Value = 62501750
$ws.Range("J99").Value = 1636.8334
$ws.Range("K99").Value = 62501750
$ws.Range("L99").Value = 1636.8334
$ws.Range("M99").Value = -62500252
$ws.Range("N99").Value = -4632.8334

$ws.Range("H134").Value = 19104.139
$ws.Range("I134").Value = 1527.8
$ws.Range("J134").Value = 79945.30499999999
$ws.Range("K134").Value = 4583.4
$ws.Range("L134").Value = 239835.915
$ws.Range("M134").Value = -2048.4
$ws.Range("N134").Value = -244905.915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3015.2827
$ws.Range("I58").Value = 583.30554
$ws.Range("K58").Value = 583.30554
$ws.Range("M58").Value = -380.30554

$ws.Range("H99").Value = 2585.2903
$ws.Range("I99").Value = 2250
$ws.Range("J99").Value = 3116.1667
$ws.Range("K99").Value = 2250
$ws.Range("L99").Value = 3116.1667
$ws.Range("M99").Value = -752
$ws.Range("N99").Value = -6112.1667

$ws.Range("H122").Value = 1013.8182
$ws.Range("I122").Value = 884.8
$ws.Range("J122").Value = 1121.3334
$ws.Range("K122").Value = 2654.4
$ws.Range("L122").Value = 3364.0002
$ws.Range("M122").Value = -204.3999999999996
$ws.Range("N122").Value = -8264.0002

$ws.Range("H126").Value = 2585.2903
$ws.Range("I126").Value = 2250
$ws.Range("J126").Value = 3116.1667
$ws.Range("K126").Value = 6750
$ws.Range("L126").Value = 9348.500100000001
$ws.Range("M126").Value = -4280
$ws.Range("N126").Value = -14288.5001

$ws.Range("H132").Value = 1344
$ws.Range("I132").Value = 811.9729599999999
$ws.Range("J132").Value = 3133.5454
$ws.Range("K132").Value = 2435.91888
$ws.Range("L132").Value = 9400.636200000001
$ws.Range("M132").Value = 94.08112000000028
$ws.Range("N132").Value = -14460.6362

$ws.Range("H134").Value = 1038.1111
$ws.Range("I134").Value = 891.7083
$ws.Range("J134").Value = 1506.6
$ws.Range("K134").Value = 2675.1249
$ws.Range("L134").Value = 4519.799999999999
$ws.Range("M134").Value = -140.1248999999998
$ws.Range("N134").Value = -9589.799999999999

$ws.Range("H136").Value = 3015.2827
$ws.Range("I136").Value = 583.30554
$ws.Range("K136").Value = 1749.91662
$ws.Range("M136").Value = 800.08338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2657.6428
$ws.Range("I132").Value = 2228.7058
$ws.Range("J132").Value = 3320.5454
$ws.Range("K132").Value = 6686.117400000001
$ws.Range("L132").Value = 9961.636200000001
$ws.Range("M132").Value = -4156.117400000001
$ws.Range("N132").Value = -15021.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1317.2667
$ws.Range("I46").Value = 1350.7273
$ws.Range("J46").Value = 1297.8948
$ws.Range("K46").Value = 1350.7273
$ws.Range("L46").Value = 1297.8948
$ws.Range("M46").Value = -1162.7273
$ws.Range("N46").Value = -1673.8948

$ws.Range("H136").Value = 2626.7795
$ws.Range("I136").Value = 1241.8269
$ws.Range("K136").Value = 3725.4807
$ws.Range("M136").Value = -1175.4807

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1328.5211
$ws.Range("I132").Value = 864.9535
$ws.Range("J132").Value = 2040.4286
$ws.Range("K132").Value = 2594.8605
$ws.Range("L132").Value = 6121.2858
$ws.Range("M132").Value = -64.86049999999977
$ws.Range("N132").Value = -11181.2858

$ws.Range("H136").Value = 2353.392
$ws.Range("I136").Value = 2327.3845
$ws.Range("J136").Value = 2437.9167
$ws.Range("K136").Value = 6982.1535
$ws.Range("L136").Value = 7313.750100000001
$ws.Range("M136").Value = -4432.1535
$ws.Range("N136").Value = -12413.7501
